$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark rows 2 and 3 as executed ("test") and update their emails/status
$ws.Range("A2").Value = "test"
$ws.Range("D2").Value = "cmcmeekan02@netvibes.com"
$ws.Range("F2").Value = "PASSED"

$ws.Range("A3").Value = "test"
$ws.Range("D3").Value = "cbreckenridge22@google.com"
$ws.Range("F3").Value = "PASSED"

# Un-mark rows 7-11 (previously "test") and set their status back to SKIPPED
$ws.Range("A7").Value = $null
$ws.Range("F7").Value = "SKIPPED"

$ws.Range("A8").Value = $null
$ws.Range("F8").Value = "SKIPPED"

$ws.Range("A9").Value = $null
$ws.Range("F9").Value = "SKIPPED"

$ws.Range("A10").Value = $null
$ws.Range("F10").Value = "SKIPPED"

$ws.Range("A11").Value = $null
$ws.Range("F11").Value = "SKIPPED"

# Move the active selection to J4 (matches the saved worksheet view)
$ws.Range("J4").Select()
